$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Espárragos" at Feria Lagunitas de
# Puerto Montt. It belongs right above the current row 52 (chronologically it
# sorts there), so insert a blank row at 52 which pushes the existing rows
# 52-80 down to 53-81, preserving all of their data and formatting.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record's data.
$ws.Range("A52").Value = 4
$ws.Range("B52").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C52").Value = "Los Lagos"
$ws.Range("D52").Value = 45236
$ws.Range("E52").Value = 10
$ws.Range("F52").Value = 300000000
$ws.Range("G52").Value = "Espárragos"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 100
$ws.Range("K52").Value = 2000
$ws.Range("L52").Value = 2000
$ws.Range("M52").Value = 2000
$ws.Range("N52").Value = "$/kilo"
$ws.Range("O52").Value = "Provincia de Linares"
$ws.Range("P52").Value = 2000
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
